$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$ws.Range("A1").Value = "Project"
$ws.Range("B1").Value = "Patch"
$ws.Range("C1").Value = "Submission"

# Data row 2
$ws.Range("A2").Value = "Cafe"
$ws.Range("B2").Value = "AstorMain-Cafe-1"
$ws.Range("C2").Value = "/Users/ruizhengu/Experiments/APR-as-AAT/anonymised-submissions/95"

# Data row 3
$ws.Range("A3").Value = "Cafe"
$ws.Range("B3").Value = "AstorMain-Cafe-2"
$ws.Range("C3").Value = "/Users/ruizhengu/Experiments/APR-as-AAT/anonymised-submissions/300"

# Build the header style once on a single cell (keeps the style table compact:
# one combined cellXf instead of one per intermediate property assignment),
# then propagate that exact style to the rest of the header row via copy/paste
# of formats only, so every header cell shares the same style index.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

$a1.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
